$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the region that is being restructured (rows 13-24, columns A-C).
# This removes all stale cells (including the entire former row 24),
# and we rebuild rows 13-23 from scratch below.
$ws.Range("A13:C24").Clear()

# Row 13
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = "Short syllabus:"
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = "Avaliação:"

# Row 18
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C18").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de laboratório, aulas de exercícios."
$ws.Range("C19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de laboratório, aulas de exercícios."
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "MF = (0,5*P1 + 0,5*P2), onde P1 e P2 são provas."
$ws.Range("C20").Value = "MF = (0,5*P1 + 0,5*P2), onde P1 e P2 são provas."
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Range("C21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Range("A22").Value = "Requisitos:"

# Row 23
$ws.Range("B23").Value = "LOQ4209 -  Engenharia da Qualidade  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4209 -  Engenharia da Qualidade  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30
